# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple worksheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 232.5
$ws.Range("I2").Value = 252
$ws.Range("J2").Value = 213
$ws.Range("K2").Value = 252
$ws.Range("L2").Value = 213
$ws.Range("M2").Value = -139
$ws.Range("N2").Value = -439
# Row 28
$ws.Range("H28").Value = 22418.5
$ws.Range("I28").Value = 1052.5
$ws.Range("J28").Value = 33101.5
$ws.Range("K28").Value = 1052.5
$ws.Range("L28").Value = 33101.5
$ws.Range("M28").Value = -567.5
$ws.Range("N28").Value = -34071.5
# Row 33
$ws.Range("H33").Value = 797.0345
$ws.Range("I33").Value = 694.5
$ws.Range("J33").Value = 964.8182
$ws.Range("K33").Value = 694.5
$ws.Range("L33").Value = 964.8182
$ws.Range("M33").Value = -465.5
$ws.Range("N33").Value = -1422.8182
# Row 62
$ws.Range("H62").Value = 2546.25
$ws.Range("I62").Value = 1650.7142
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 1650.7142
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -1026.7142
$ws.Range("N62").Value = -5048
# Row 65
$ws.Range("H65").Value = 2546.25
$ws.Range("I65").Value = 1650.7142
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 8253.571
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -5133.571
$ws.Range("N65").Value = -25240
# Row 98
$ws.Range("H98").Value = 1960
$ws.Range("I98").Value = 1600
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 1600
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = -102
$ws.Range("N98").Value = -5496
# Row 107
$ws.Range("H107").Value = 1488.4
$ws.Range("I107").Value = 1615.8182
$ws.Range("J107").Value = 1138
$ws.Range("K107").Value = 1615.8182
$ws.Range("L107").Value = 1138
$ws.Range("M107").Value = 304.1818000000001
$ws.Range("N107").Value = -4978
# Row 122
$ws.Range("H122").Value = 1960
$ws.Range("I122").Value = 1600
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4800
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -2350
$ws.Range("N122").Value = -12400
# Row 138
$ws.Range("H138").Value = 183147.55
$ws.Range("J138").Value = 280614.22
$ws.Range("L138").Value = 841842.6599999999
$ws.Range("N138").Value = -852122.6599999999

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3404
$ws.Range("I61").Value = 3404
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3404
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3192
$ws.Range("N61").ClearContents()
# Row 74
$ws.Range("H74").Value = 1620.8422
$ws.Range("I74").Value = 1516.5883
$ws.Range("J74").Value = 2507
$ws.Range("K74").Value = 1516.5883
$ws.Range("L74").Value = 2507
$ws.Range("M74").Value = -642.5882999999999
$ws.Range("N74").Value = -4255
# Row 77
$ws.Range("H77").Value = 1620.8422
$ws.Range("I77").Value = 1516.5883
$ws.Range("J77").Value = 2507
$ws.Range("K77").Value = 7582.941499999999
$ws.Range("L77").Value = 12535
$ws.Range("M77").Value = -3214.941499999999
$ws.Range("N77").Value = -21271
# Row 136
$ws.Range("H136").Value = 3404
$ws.Range("I136").Value = 3404
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10212
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7662
$ws.Range("N136").Value = 0

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1579.4445
$ws.Range("I107").Value = 863
$ws.Range("K107").Value = 863
$ws.Range("M107").Value = 1057

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 37500.75
$ws.Range("I4").Value = 5001
$ws.Range("J4").Value = 48334
$ws.Range("K4").Value = 5001
$ws.Range("L4").Value = 48334
$ws.Range("M4").Value = -4889
$ws.Range("N4").Value = -48558
# Row 31
$ws.Range("H31").Value = 1529.6364
$ws.Range("I31").Value = 916.43335
$ws.Range("J31").Value = 2843.6428
$ws.Range("K31").Value = 916.43335
$ws.Range("L31").Value = 2843.6428
$ws.Range("M31").Value = -621.43335
$ws.Range("N31").Value = -3433.6428
# Row 34
$ws.Range("H34").Value = 1529.6364
$ws.Range("I34").Value = 916.43335
$ws.Range("J34").Value = 2843.6428
$ws.Range("K34").Value = 916.43335
$ws.Range("L34").Value = 2843.6428
$ws.Range("M34").Value = -714.43335
$ws.Range("N34").Value = -3247.6428
# Row 58
$ws.Range("H58").Value = 1362.4445
$ws.Range("I58").Value = 912.5
$ws.Range("J58").Value = 1491
$ws.Range("K58").Value = 912.5
$ws.Range("L58").Value = 1491
$ws.Range("M58").Value = -709.5
$ws.Range("N58").Value = -1897
# Row 107
$ws.Range("H107").Value = 478.2258
$ws.Range("I107").Value = 346.5909
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 346.5909
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1573.4091
$ws.Range("N107").Value = -4640
# Row 132
$ws.Range("H132").Value = 11113286
$ws.Range("I132").Value = 1131.6666
$ws.Range("J132").Value = 18521388
$ws.Range("K132").Value = 3394.9998
$ws.Range("L132").Value = 55564164
$ws.Range("M132").Value = -864.9998000000001
$ws.Range("N132").Value = -55569224
# Row 136
$ws.Range("H136").Value = 1362.4445
$ws.Range("I136").Value = 912.5
$ws.Range("J136").Value = 1491
$ws.Range("K136").Value = 2737.5
$ws.Range("L136").Value = 4473
$ws.Range("M136").Value = -187.5
$ws.Range("N136").Value = -9573

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1238.8
$ws.Range("I5").Value = 1238.8
$ws.Range("K5").Value = 3716.4
$ws.Range("M5").Value = -3604.4
# Row 135
$ws.Range("H135").Value = 1238.8
$ws.Range("I135").Value = 1238.8
$ws.Range("K135").Value = 11149.2
$ws.Range("M135").Value = -8614.199999999999
# Row 138
$ws.Range("H138").Value = 3847.6316
$ws.Range("I138").Value = 1244
$ws.Range("J138").Value = 4242.121
$ws.Range("K138").Value = 3732
$ws.Range("L138").Value = 12726.363
$ws.Range("M138").Value = 1408
$ws.Range("N138").Value = -23006.363

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 23778.223
$ws.Range("J2").Value = 47251
$ws.Range("L2").Value = 47251
$ws.Range("N2").Value = -47475
# Row 22
$ws.Range("H22").Value = 4511.967
$ws.Range("I22").Value = 792.8570999999999
$ws.Range("J22").Value = 7766.1875
$ws.Range("K22").Value = 792.8570999999999
$ws.Range("L22").Value = 7766.1875
$ws.Range("M22").Value = -497.8570999999999
$ws.Range("N22").Value = -8356.1875
# Row 27
$ws.Range("H27").Value = 4511.967
$ws.Range("I27").Value = 792.8570999999999
$ws.Range("J27").Value = 7766.1875
$ws.Range("K27").Value = 792.8570999999999
$ws.Range("L27").Value = 7766.1875
$ws.Range("M27").Value = -685.8570999999999
$ws.Range("N27").Value = -7980.1875
# Row 93
$ws.Range("H93").Value = 6136.25
$ws.Range("I93").Value = 7922.125
$ws.Range("J93").Value = 2564.5
$ws.Range("K93").Value = 7922.125
$ws.Range("L93").Value = 2564.5
$ws.Range("M93").Value = -6674.125
$ws.Range("N93").Value = -5060.5
# Row 136
$ws.Range("H136").Value = 2777.9167
$ws.Range("I136").Value = 3919.1667
$ws.Range("J136").Value = 1636.6666
$ws.Range("K136").Value = 11757.5001
$ws.Range("L136").Value = 4909.9998
$ws.Range("M136").Value = -9207.500100000001
$ws.Range("N136").Value = -10009.9998

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 3243.543
$ws.Range("I136").Value = 2938.1875
$ws.Range("J136").Value = 3500.6843
$ws.Range("K136").Value = 8814.5625
$ws.Range("L136").Value = 10502.0529
$ws.Range("M136").Value = -6264.5625
$ws.Range("N136").Value = -15602.0529
